$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "dsds"
$ws.Range("L2").ClearContents()
$ws.Range("X2").ClearContents()

# --- Row 3 ---
$ws.Range("N3").ClearContents()
$ws.Range("V3").Value = "kkkk"
$ws.Range("Y3").Value = "ds"
$ws.Range("AC3").Value = "hhh"

# --- Row 4 ---
$ws.Range("B4").Value = "50.50.2020"
$ws.Range("D4").Value = "fa"
$ws.Range("L4").Value = "asd"
$ws.Range("AB4").Value = "aaa"

# --- Selection ---
$ws.Range("B5").Select() | Out-Null
